$d = $word.ActiveDocument

# Update the date line at the top of the document
$d.Content.Find.Execute("2024-04-30 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-05-01 Wednesday", 2)

# Update each multiplication expression cell in the (only) table, addressed by
# row/column so that values which collide with other cells old/new text are not
# mismatched by a global find/replace.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "13×30="
$t.Cell(1, 2).Range.Text = "80×91="
$t.Cell(1, 3).Range.Text = "68×37="
$t.Cell(1, 4).Range.Text = "54×77="
$t.Cell(1, 5).Range.Text = "30×51="
$t.Cell(5, 1).Range.Text = "78×55="
$t.Cell(5, 2).Range.Text = "11×41="
$t.Cell(5, 3).Range.Text = "18×37="
$t.Cell(5, 4).Range.Text = "40×21="
$t.Cell(5, 5).Range.Text = "77×41="
$t.Cell(10, 1).Range.Text = "78×23="
$t.Cell(10, 2).Range.Text = "29×96="
$t.Cell(10, 3).Range.Text = "69×60="
$t.Cell(10, 4).Range.Text = "46×53="
$t.Cell(10, 5).Range.Text = "52×87="
$t.Cell(15, 1).Range.Text = "95×46="
$t.Cell(15, 2).Range.Text = "93×83="
$t.Cell(15, 3).Range.Text = "38×66="
$t.Cell(15, 4).Range.Text = "87×85="
$t.Cell(15, 5).Range.Text = "39×22="
$t.Cell(20, 1).Range.Text = "99×59="
$t.Cell(20, 2).Range.Text = "82×18="
$t.Cell(20, 3).Range.Text = "76×86="
$t.Cell(20, 4).Range.Text = "26×44="
$t.Cell(20, 5).Range.Text = "66×47="
